# Apply the edits described by the commit:
#  - B5 / B6 on Sheet1: "确定用例" -> "确定手机端用例"
#  - Sheet1 selection moves from M13 to H12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two shared-string cells that mention "确定用例" -> "确定手机端用例"
$b5 = $ws.Range("B5").Value2
if ($b5 -like "*确定用例*") {
    $ws.Range("B5").Value = $b5 -replace "确定用例", "确定手机端用例"
}

$b6 = $ws.Range("B6").Value2
if ($b6 -like "*确定用例*") {
    $ws.Range("B6").Value = $b6 -replace "确定用例", "确定手机端用例"
}

# Move the active selection on Sheet1 from M13 to H12
[void]$ws.Activate()
[void]$ws.Range("H12").Select()
